$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Cells.Item(19, 2).Value = 6884522
$ws.Cells.Item(19, 5).Value = "Schalke II"
$ws.Cells.Item(19, 6).Value = "Wuppertaler"
$ws.Cells.Item(19, 7).Value = 1
$ws.Cells.Item(19, 8).Value = 3
$ws.Cells.Item(19, 9).Value = 1
$ws.Cells.Item(19, 10).Value = 1
$ws.Cells.Item(19, 11).Value = "A"
$ws.Cells.Item(19, 12).Value = 2.6
$ws.Cells.Item(19, 13).Value = 3.6
$ws.Cells.Item(19, 14).Value = 2.2
$ws.Cells.Item(19, 15).Value = 2.9
$ws.Cells.Item(19, 16).Value = 3.6
$ws.Cells.Item(19, 17).Value = 2
$ws.Cells.Item(19, 18).Value = 0.25
$ws.Cells.Item(19, 19).Value = 1.975
$ws.Cells.Item(19, 20).Value = 1.825
$ws.Cells.Item(19, 21).Value = 3.25
$ws.Cells.Item(19, 22).Value = 1.8
$ws.Cells.Item(19, 23).Value = 2
$ws.Cells.Item(19, 24).Value = -1
$ws.Cells.Item(19, 25).Value = -1
$ws.Cells.Item(19, 26).Value = 1
$ws.Cells.Item(19, 27).Value = -1
$ws.Cells.Item(19, 28).Value = 0.825
$ws.Cells.Item(19, 29).Value = 0.8
$ws.Cells.Item(19, 30).Value = -1

# Row 20
$ws.Cells.Item(20, 2).Value = 6884488
$ws.Cells.Item(20, 5).Value = "Alemannia Aachen"
$ws.Cells.Item(20, 6).Value = "Borussia Mgladbach II"
$ws.Cells.Item(20, 7).Value = 2
$ws.Cells.Item(20, 8).Value = 2
$ws.Cells.Item(20, 9).Value = 2
$ws.Cells.Item(20, 10).Value = 1
$ws.Cells.Item(20, 11).Value = "D"
$ws.Cells.Item(20, 12).Value = 2.4
$ws.Cells.Item(20, 13).Value = 3.6
$ws.Cells.Item(20, 14).Value = 2.4
$ws.Cells.Item(20, 15).Value = 1.8
$ws.Cells.Item(20, 16).Value = 3.6
$ws.Cells.Item(20, 17).Value = 3.4
$ws.Cells.Item(20, 18).Value = -0.5
$ws.Cells.Item(20, 19).Value = 1.9
$ws.Cells.Item(20, 20).Value = 1.95
$ws.Cells.Item(20, 21).Value = 3
$ws.Cells.Item(20, 22).Value = 2.025
$ws.Cells.Item(20, 23).Value = 1.825
$ws.Cells.Item(20, 24).Value = -1
$ws.Cells.Item(20, 25).Value = 2.6
$ws.Cells.Item(20, 26).Value = -1
$ws.Cells.Item(20, 27).Value = -1
$ws.Cells.Item(20, 28).Value = 0.95
$ws.Cells.Item(20, 29).Value = 1.025
$ws.Cells.Item(20, 30).Value = -1

# Row 21
$ws.Cells.Item(21, 2).Value = 6886885
$ws.Cells.Item(21, 5).Value = "SC Wiedenbruck"
$ws.Cells.Item(21, 6).Value = "Fortuna Dusseldorf II"
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(21, 8).Value = 4
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 2
$ws.Cells.Item(21, 11).Value = "A"
$ws.Cells.Item(21, 12).Value = 1.95
$ws.Cells.Item(21, 13).Value = 3.5
$ws.Cells.Item(21, 14).Value = 3.2
$ws.Cells.Item(21, 15).Value = 1.727
$ws.Cells.Item(21, 16).Value = 3.6
$ws.Cells.Item(21, 17).Value = 3.8
$ws.Cells.Item(21, 18).Value = -0.75
$ws.Cells.Item(21, 19).Value = 2
$ws.Cells.Item(21, 20).Value = 1.8
$ws.Cells.Item(21, 21).Value = 3
$ws.Cells.Item(21, 22).Value = 1.95
$ws.Cells.Item(21, 23).Value = 1.85
$ws.Cells.Item(21, 24).Value = -1
$ws.Cells.Item(21, 25).Value = -1
$ws.Cells.Item(21, 26).Value = 2.8
$ws.Cells.Item(21, 27).Value = -1
$ws.Cells.Item(21, 28).Value = 0.8
$ws.Cells.Item(21, 29).Value = 0.95
$ws.Cells.Item(21, 30).Value = -1

# Row 22
$ws.Cells.Item(22, 2).Value = 6886888
$ws.Cells.Item(22, 5).Value = "SV Rodinghausen"
$ws.Cells.Item(22, 6).Value = "Duren"
$ws.Cells.Item(22, 7).Value = 2
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 1
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = "H"
$ws.Cells.Item(22, 12).Value = 1.666
$ws.Cells.Item(22, 13).Value = 3.8
$ws.Cells.Item(22, 14).Value = 4
$ws.Cells.Item(22, 15).Value = 1.666
$ws.Cells.Item(22, 16).Value = 3.75
$ws.Cells.Item(22, 17).Value = 4
$ws.Cells.Item(22, 18).Value = -0.75
$ws.Cells.Item(22, 19).Value = 1.925
$ws.Cells.Item(22, 20).Value = 1.875
$ws.Cells.Item(22, 21).Value = 2.75
$ws.Cells.Item(22, 22).Value = 1.825
$ws.Cells.Item(22, 23).Value = 1.975
$ws.Cells.Item(22, 24).Value = 0.6659999999999999
$ws.Cells.Item(22, 25).Value = -1
$ws.Cells.Item(22, 26).Value = -1
$ws.Cells.Item(22, 27).Value = 0.925
$ws.Cells.Item(22, 28).Value = -1
$ws.Cells.Item(22, 29).Value = -1
$ws.Cells.Item(22, 30).Value = 0.9750000000000001

# Row 297
$ws.Cells.Item(297, 2).Value = 7088567
$ws.Cells.Item(297, 5).Value = "Alemannia Aachen"
$ws.Cells.Item(297, 6).Value = "SSVg Velbert"
$ws.Cells.Item(297, 7).Value = 2
$ws.Cells.Item(297, 8).Value = 0
$ws.Cells.Item(297, 9).Value = 1
$ws.Cells.Item(297, 10).Value = 0
$ws.Cells.Item(297, 11).Value = "H"
$ws.Cells.Item(297, 12).Value = 1.2
$ws.Cells.Item(297, 13).Value = 6
$ws.Cells.Item(297, 14).Value = 9
$ws.Cells.Item(297, 15).Value = 1.333
$ws.Cells.Item(297, 16).Value = 5.25
$ws.Cells.Item(297, 17).Value = 6
$ws.Cells.Item(297, 18).Value = -1.5
$ws.Cells.Item(297, 19).Value = 1.9
$ws.Cells.Item(297, 20).Value = 1.9
$ws.Cells.Item(297, 21).Value = 3.5
$ws.Cells.Item(297, 22).Value = 2
$ws.Cells.Item(297, 23).Value = 1.8
$ws.Cells.Item(297, 24).Value = 0.333
$ws.Cells.Item(297, 25).Value = -1
$ws.Cells.Item(297, 26).Value = -1
$ws.Cells.Item(297, 27).Value = 0.8999999999999999
$ws.Cells.Item(297, 28).Value = -1
$ws.Cells.Item(297, 29).Value = -1
$ws.Cells.Item(297, 30).Value = 0.8

# Row 298
$ws.Cells.Item(298, 2).Value = 7088568
$ws.Cells.Item(298, 5).Value = "SC Paderborn 07 II"
$ws.Cells.Item(298, 6).Value = "RotWeiss Oberhausen"
$ws.Cells.Item(298, 7).Value = 1
$ws.Cells.Item(298, 8).Value = 1
$ws.Cells.Item(298, 9).Value = 0
$ws.Cells.Item(298, 10).Value = 0
$ws.Cells.Item(298, 11).Value = "D"
$ws.Cells.Item(298, 12).Value = 2.3
$ws.Cells.Item(298, 13).Value = 3.75
$ws.Cells.Item(298, 14).Value = 2.4
$ws.Cells.Item(298, 15).Value = 2.375
$ws.Cells.Item(298, 16).Value = 3.8
$ws.Cells.Item(298, 17).Value = 2.3
$ws.Cells.Item(298, 18).Value = 0
$ws.Cells.Item(298, 19).Value = 1.95
$ws.Cells.Item(298, 20).Value = 1.9
$ws.Cells.Item(298, 21).Value = 3
$ws.Cells.Item(298, 22).Value = 1.825
$ws.Cells.Item(298, 23).Value = 2.025
$ws.Cells.Item(298, 24).Value = -1
$ws.Cells.Item(298, 25).Value = 2.8
$ws.Cells.Item(298, 26).Value = -1
$ws.Cells.Item(298, 27).Value = 0
$ws.Cells.Item(298, 28).Value = 0
$ws.Cells.Item(298, 29).Value = -1
$ws.Cells.Item(298, 30).Value = 1.025

# Row 299
$ws.Cells.Item(299, 2).Value = 7088569
$ws.Cells.Item(299, 5).Value = "SC Wiedenbruck"
$ws.Cells.Item(299, 6).Value = "SV Lippstadt 08"
$ws.Cells.Item(299, 7).Value = 1
$ws.Cells.Item(299, 8).Value = 1
$ws.Cells.Item(299, 9).Value = 1
$ws.Cells.Item(299, 10).Value = 1
$ws.Cells.Item(299, 11).Value = "D"
$ws.Cells.Item(299, 12).Value = 1.571
$ws.Cells.Item(299, 13).Value = 4
$ws.Cells.Item(299, 14).Value = 4.5
$ws.Cells.Item(299, 15).Value = 1.8
$ws.Cells.Item(299, 16).Value = 3.9
$ws.Cells.Item(299, 17).Value = 3.4
$ws.Cells.Item(299, 18).Value = -0.5
$ws.Cells.Item(299, 19).Value = 1.85
$ws.Cells.Item(299, 20).Value = 2
$ws.Cells.Item(299, 21).Value = 3.25
$ws.Cells.Item(299, 22).Value = 1.85
$ws.Cells.Item(299, 23).Value = 2
$ws.Cells.Item(299, 24).Value = -1
$ws.Cells.Item(299, 25).Value = 2.9
$ws.Cells.Item(299, 26).Value = -1
$ws.Cells.Item(299, 27).Value = -1
$ws.Cells.Item(299, 28).Value = 1
$ws.Cells.Item(299, 29).Value = -1
$ws.Cells.Item(299, 30).Value = 1

# Row 300
$ws.Cells.Item(300, 2).Value = 7091975
$ws.Cells.Item(300, 5).Value = "Wuppertaler"
$ws.Cells.Item(300, 6).Value = "Duren"
$ws.Cells.Item(300, 7).Value = 1
$ws.Cells.Item(300, 8).Value = 1
$ws.Cells.Item(300, 9).Value = 0
$ws.Cells.Item(300, 10).Value = 1
$ws.Cells.Item(300, 11).Value = "D"
$ws.Cells.Item(300, 12).Value = 1.6
$ws.Cells.Item(300, 13).Value = 4
$ws.Cells.Item(300, 14).Value = 4.2
$ws.Cells.Item(300, 15).Value = 1.7
$ws.Cells.Item(300, 16).Value = 4
$ws.Cells.Item(300, 17).Value = 3.8
$ws.Cells.Item(300, 18).Value = -0.75
$ws.Cells.Item(300, 19).Value = 1.9
$ws.Cells.Item(300, 20).Value = 1.9
$ws.Cells.Item(300, 21).Value = 3.25
$ws.Cells.Item(300, 22).Value = 1.825
$ws.Cells.Item(300, 23).Value = 1.975
$ws.Cells.Item(300, 24).Value = -1
$ws.Cells.Item(300, 25).Value = 3
$ws.Cells.Item(300, 26).Value = -1
$ws.Cells.Item(300, 27).Value = -1
$ws.Cells.Item(300, 28).Value = 0.8999999999999999
$ws.Cells.Item(300, 29).Value = -1
$ws.Cells.Item(300, 30).Value = 0.9750000000000001

# Row 301
$ws.Cells.Item(301, 2).Value = 7091976
$ws.Cells.Item(301, 5).Value = "Borussia Mgladbach II"
$ws.Cells.Item(301, 6).Value = "1 FC Bocholt"
$ws.Cells.Item(301, 7).Value = 3
$ws.Cells.Item(301, 8).Value = 3
$ws.Cells.Item(301, 9).Value = 1
$ws.Cells.Item(301, 10).Value = 1
$ws.Cells.Item(301, 11).Value = "D"
$ws.Cells.Item(301, 12).Value = 2.6
$ws.Cells.Item(301, 13).Value = 4
$ws.Cells.Item(301, 14).Value = 2.1
$ws.Cells.Item(301, 15).Value = 2.5
$ws.Cells.Item(301, 16).Value = 4.1
$ws.Cells.Item(301, 17).Value = 2.15
$ws.Cells.Item(301, 18).Value = 0.25
$ws.Cells.Item(301, 19).Value = 1.8
$ws.Cells.Item(301, 20).Value = 2
$ws.Cells.Item(301, 21).Value = 3.25
$ws.Cells.Item(301, 22).Value = 1.95
$ws.Cells.Item(301, 23).Value = 1.85
$ws.Cells.Item(301, 24).Value = -1
$ws.Cells.Item(301, 25).Value = 3.1
$ws.Cells.Item(301, 26).Value = -1
$ws.Cells.Item(301, 27).Value = 0.4
$ws.Cells.Item(301, 28).Value = -0.5
$ws.Cells.Item(301, 29).Value = 0.95
$ws.Cells.Item(301, 30).Value = -1

# Row 302
$ws.Cells.Item(302, 2).Value = 7091978
$ws.Cells.Item(302, 5).Value = "SV Rodinghausen"
$ws.Cells.Item(302, 6).Value = "SC Fortuna Kln"
$ws.Cells.Item(302, 7).Value = 1
$ws.Cells.Item(302, 8).Value = 0
$ws.Cells.Item(302, 9).Value = 1
$ws.Cells.Item(302, 10).Value = 0
$ws.Cells.Item(302, 11).Value = "H"
$ws.Cells.Item(302, 12).Value = 2.25
$ws.Cells.Item(302, 13).Value = 4
$ws.Cells.Item(302, 14).Value = 2.4
$ws.Cells.Item(302, 15).Value = 1.95
$ws.Cells.Item(302, 16).Value = 3.9
$ws.Cells.Item(302, 17).Value = 2.9
$ws.Cells.Item(302, 18).Value = -0.25
$ws.Cells.Item(302, 19).Value = 1.825
$ws.Cells.Item(302, 20).Value = 2.025
$ws.Cells.Item(302, 21).Value = 3.25
$ws.Cells.Item(302, 22).Value = 1.85
$ws.Cells.Item(302, 23).Value = 2
$ws.Cells.Item(302, 24).Value = 0.95
$ws.Cells.Item(302, 25).Value = -1
$ws.Cells.Item(302, 26).Value = -1
$ws.Cells.Item(302, 27).Value = 0.825
$ws.Cells.Item(302, 28).Value = -1
$ws.Cells.Item(302, 29).Value = -1
$ws.Cells.Item(302, 30).Value = 1

# Row 303
$ws.Cells.Item(303, 2).Value = 7091980
$ws.Cells.Item(303, 5).Value = "FC WegbergBeeck"
$ws.Cells.Item(303, 6).Value = "Gutersloh 2000"
$ws.Cells.Item(303, 7).Value = 1
$ws.Cells.Item(303, 8).Value = 2
$ws.Cells.Item(303, 9).Value = 0
$ws.Cells.Item(303, 10).Value = 1
$ws.Cells.Item(303, 11).Value = "A"
$ws.Cells.Item(303, 12).Value = 2.3
$ws.Cells.Item(303, 13).Value = 3.75
$ws.Cells.Item(303, 14).Value = 2.4
$ws.Cells.Item(303, 15).Value = 2.9
$ws.Cells.Item(303, 16).Value = 4
$ws.Cells.Item(303, 17).Value = 1.85
$ws.Cells.Item(303, 18).Value = 0.5
$ws.Cells.Item(303, 19).Value = 1.85
$ws.Cells.Item(303, 20).Value = 1.95
$ws.Cells.Item(303, 21).Value = 3.25
$ws.Cells.Item(303, 22).Value = 1.9
$ws.Cells.Item(303, 23).Value = 1.9
$ws.Cells.Item(303, 24).Value = -1
$ws.Cells.Item(303, 25).Value = -1
$ws.Cells.Item(303, 26).Value = 0.8500000000000001
$ws.Cells.Item(303, 27).Value = -1
$ws.Cells.Item(303, 28).Value = 0.95
$ws.Cells.Item(303, 29).Value = -0.5
$ws.Cells.Item(303, 30).Value = 0.45
